# "Fix Tax & Management fee"
# Insert two new columns (IVA %, Importe Neto) before the existing
# "Management Fee" columns, and rework the Management Fee formula so it is
# computed on the net amount (after VAT) instead of the gross amount.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert two blank columns at L:M -------------------------------
# Everything that used to live in L..R shifts right to N..T.
$ws.Columns("L:M").Insert()

# --- 2. Header row (row 2) for the two new columns ---------------------
$ws.Range("L2").Value = "IVA"
$ws.Range("M2").Value = "Importe Neto"

# Carry over the header style used by the rest of row 2.
$ws.Range("L2").Style = $ws.Range("K2").Style
$ws.Range("M2").Style = $ws.Range("K2").Style

# --- 3. Row 1 (subtotal row) formatting for the new columns -------------
# New L1/M1 (and the shifted-in N1) reuse the same formatting the old L1
# (now N1) already had.
$ws.Range("L1").Style = $ws.Range("N1").Style
$ws.Range("M1").Style = $ws.Range("N1").Style

# --- 4. Row 3 (data/formula row) ---------------------------------------
# L3: IVA % input cell -> new "Percent" style (creates the Porcentaje
# cell style + links it to the xf, as in the target workbook).
$ws.Range("L3").Value = $null
$ws.Range("L3").Style = "Percent"
$ws.Range("L3").HorizontalAlignment = -4152  # xlRight
$ws.Range("L3").VerticalAlignment = -4160    # xlTop
$ws.Range("L3").Borders.LineStyle = $ws.Range("N3").Borders.LineStyle

# M3: Importe Neto = Importe / (1 + IVA%)
$ws.Range("M3").Formula = "=J3/(1+L3)"
$ws.Range("M3").Style = $ws.Range("J3").Style

# N3: Management Fee % input cell - keeps the formatting the old L3 had.
$ws.Range("N3").Style = $ws.Range("L3").Style
$ws.Range("N3").NumberFormat = "0%"
$ws.Range("N3").Borders.LineStyle = 1
$ws.Range("N3").Value = $null

# O3: Management Fee amount = Importe Neto * Management Fee %
$ws.Range("O3").Formula = "=M3*N3"
$ws.Range("O3").Style = $ws.Range("J3").Style

# --- 5. Subtotal formula for the Management Fee column (row 1) ---------
# The SUBTOTAL formula that used to sit on M1 now belongs on O1 (it moved
# automatically with the column insert), make sure it is correct.
$ws.Range("O1").Formula = "=SUBTOTAL(9,O3:O99998)"

# --- 6. Re-point the AutoFilter / _FilterDatabase to the new extent ----
$ws.AutoFilterMode = $false
$ws.Range("A2:T2").AutoFilter()

$fd = $wb.Names.Item("Ingresos!_FilterDatabase")
$fd.RefersTo = "=Ingresos!`$A`$2:`$T`$2"
